$wb = $excel.ActiveWorkbook

$violSheet = $wb.Worksheets.Item("3 - Table Block - Violations")
[void]$violSheet.Range("B1:B13").Select()

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "3-TableBlock-TableMetricIdCol"

$newSheet.Range("B1").Value = "3.4. - TABLE_METRIC_ID_COL"
$newSheet.Range("B2").Value = "* Block Name = TABLE_METRIC_ID_COL"
$newSheet.Range("B3").Value = "* Options :"
$newSheet.Range("B4").Value = "- QID : list of Quality indicators (BC or TC or RULE) separated by |"
$newSheet.Range("B5").Value = "- SID : list of Sizing measures separated by |"
$newSheet.Range("B6").Value = "- BID : list of Background facts separated by |"
$newSheet.Range("B7").Value = "- LEVEL : can be APPLICATION or MODULES or TECHNOLOGIES (by default APPLICATION if option not present)"
$newSheet.Range("B9").Value = "- VARIATION = VALUE or PERCENT or BOTH (PERCENT by default)"
$newSheet.Range("B8").Value = "- SNAPSHOT=CURRENT (only CURRENT SNAPSHOT) or PREVIOUS (only PREVIOUS SNAPSHOT) or BOTH (CURRENT and PREVIOUS SNAPSHOT, default option)"
$newSheet.Range("B10").Value = "- HEADER=SHORT, SHORT name is taken if exists, name otherwise"
$newSheet.Range("B12").Value = "RepGen:TABLE;TABLE_METRIC_ID_COL;QID=60017|60014,SID=10151|67010,BID=66061,LEVEL=APPLICATION,SNAPSHOT=BOTH,VARIATION=BOTH,HEADER=SHORT"

$srcSheet = $wb.Worksheets.Item("2 - Text blocks")

$srcSheet.Range("B2").Copy()
$newSheet.Range("B1").PasteSpecial(-4122)

$srcSheet.Range("B3").Copy()
$newSheet.Range("B2:B3").PasteSpecial(-4122)
$newSheet.Range("B11:B12").PasteSpecial(-4122)

$srcSheet.Range("B61").Copy()
$newSheet.Range("B4:B10").PasteSpecial(-4122)

$newSheet.Rows.Item(1).RowHeight = 18.75
$newSheet.PageSetup.Orientation = 1

[void]$newSheet.Range("B12").Select()

Write-Output "done"
